# Powerpoint writer: consolidate text run nodes.
# This merges each "word" run with the single-space run that immediately
# follows it (when present), reducing the number of <a:r> nodes emitted
# for the title placeholders and the "an/An image" captions.
#
# We drive this purely through TextRange.Characters(start, length).Text = ...
# which (per observed COM-interop semantics) performs a minimal in-place
# edit of the underlying runs rather than a full text-frame rebuild, so
# assigning the already-current text back (just merging run boundaries)
# causes adjacent runs covered by a single Characters() call to collapse
# into one run, exactly matching the target OOXML.

function Merge-Runs {
    param($TextRange, $Segments)
    $pos = 1
    foreach ($seg in $Segments) {
        $len = $seg.Length
        if ($len -gt 0) {
            $TextRange.Characters($pos, $len).Text = $seg
        }
        $pos = $pos + $len
    }
}

$p = $ppt.ActivePresentation

# --- Title placeholders: "Slide N (...)" on every slide ---------------
$titles = @{
    1  = @("Slide ", "1 ", "(Content)")
    2  = @("Slide ", "2 ", "(Content)")
    3  = @("Slide ", "3 ", "(Content)")
    4  = @("Slide ", "4 ", "(Content)")
    5  = @("Slide ", "5 ", "(Two ", "Content)")
    6  = @("Slide ", "6 ", "(Two ", "Content ", "Right)")
    7  = @("Slide ", "7 ", "(Content ", "with ", "Caption)")
    8  = @("Slide ", "8 ", "(Comparison)")
    9  = @("Slide ", "10 ", "(Content)")
    10 = @("Slide ", "11 ", "(Content)")
    11 = @("Slide ", "12 ", "(Content)")
}

function Find-ShapeByName {
    param($Slide, $Name)
    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $sh = $Slide.Shapes.Item($i)
        if ($sh.Name -eq $Name) {
            return $sh
        }
    }
    return $null
}

foreach ($idx in 1..11) {
    $slide = $p.Slides.Item($idx)
    $titleShape = Find-ShapeByName $slide "Title 1"
    $tr = $titleShape.TextFrame.TextRange
    Merge-Runs $tr $titles[$idx]
}

# --- "an image" / "An image" captions ----------------------------------
$captions = @(
    @{ Slide = 6; Segments = @("an ", "image") }
    @{ Slide = 7; Segments = @("An ", "image") }
    @{ Slide = 8; Segments = @("An ", "image") }
)

foreach ($cap in $captions) {
    $slide = $p.Slides.Item($cap.Slide)
    $shape = Find-ShapeByName $slide "TextBox 3"
    $tr = $shape.TextFrame.TextRange
    Merge-Runs $tr $cap.Segments
}
